$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.472
$ws.Range("D9").Value = -8.317000000000002
$ws.Range("B11").Value = 5.681999999999999
$ws.Range("C11").Value = -12.064
$ws.Range("B12").Value = 4.951000000000001
$ws.Range("D13").Value = -8.045
$ws.Range("D14").Value = -8.196999999999999
$ws.Range("B15").Value = 5.198
$ws.Range("D19").Value = -8.093999999999999
$ws.Range("D21").Value = -8.404
$ws.Range("D22").Value = -8.215
$ws.Range("C23").Value = -12.695
$ws.Range("D24").Value = -7.339999999999999
$ws.Range("D26").Value = -7.358999999999999
$ws.Range("B27").Value = 5.378
$ws.Range("B28").Value = 5.684
$ws.Range("C28").Value = -12.607
$ws.Range("B31").Value = 5.487
$ws.Range("B32").Value = 5.875999999999999
$ws.Range("C32").Value = -11.925
$ws.Range("C34").Value = -11.849
$ws.Range("B36").Value = 9.260999999999999
$ws.Range("C36").Value = -12.869
$ws.Range("C37").Value = -12.912
$ws.Range("B38").Value = 6.126
$ws.Range("D38").Value = -7.840000000000001
$ws.Range("D41").Value = -8.489000000000001
$ws.Range("C42").Value = -12.761
$ws.Range("B46").Value = 6.167
$ws.Range("C49").Value = -12.983
$ws.Range("D52").Value = -7.876
$ws.Range("B54").Value = 5.335999999999999
$ws.Range("C54").Value = -12.942
$ws.Range("B55").Value = 4.726
$ws.Range("B56").Value = 5.187
$ws.Range("D56").Value = -8.130999999999998
$ws.Range("B67").Value = 5.754
$ws.Range("B69").Value = 5.535
$ws.Range("D71").Value = -7.547
$ws.Range("B72").Value = 5.779999999999999
$ws.Range("D72").Value = -7.617
$ws.Range("B73").Value = 7.729000000000001
$ws.Range("C78").Value = -12.434
$ws.Range("D78").Value = -7.467999999999999
$ws.Range("C80").Value = -11.332
$ws.Range("B83").Value = 5.107
$ws.Range("D83").Value = -8.422999999999998
$ws.Range("D85").Value = -8.359000000000002
$ws.Range("B86").Value = 5.037
$ws.Range("D86").Value = -8.422999999999998
$ws.Range("D90").Value = -7.361999999999999
$ws.Range("B91").Value = 5.278
$ws.Range("B93").Value = 5.382
$ws.Range("D96").Value = -7.417999999999999
$ws.Range("C97").Value = -11.716
$ws.Range("B99").Value = 5.226000000000001
$ws.Range("C99").Value = -11.278
$ws.Range("C100").Value = -11.898
$ws.Range("C101").Value = -12.481
$ws.Range("D103").Value = -8.346
$ws.Range("B104").Value = 8.209
$ws.Range("B105").Value = 8.241000000000001
